{"js": "// Remove the \"This is a fancy text example...\" paragraph entirely, and\n// have the following paragraph (\"LaTeX can be written like LaTeX.\") take\n// over the style (\"FirstParagraph\") that the removed paragraph used to have.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet target = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"This is a fancy text example\") !== -1) {\n    target = items[i];\n    break;\n  }\n}\n\nif (target) {\n  const next = target.getNext();\n  next.load(\"style\");\n  await context.sync();\n\n  next.style = target.style;\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"This is a fancy text example...\" paragraph entirely, and\n# have the following paragraph (\"LaTeX can be written like LaTeX.\") take\n# over the style (\"FirstParagraph\") that the removed paragraph used to have.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*This is a fancy text example*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $next = $target.Next()\n    $next.Style = $target.Style\n    $target.Range.Delete()\n}\n"}
